$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-11-29"

# Update the November row label
$ws.Range("A12").Value = "November (through 11-29)"

# Update November row (row 12) figures
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = 76
$ws.Range("D12").Value = 105
$ws.Range("E12").Value = 67
$ws.Range("G12").Value = 205
$ws.Range("H12").Value = 192

# Update Total row (row 13) figures
$ws.Range("B13").Value = 290
$ws.Range("C13").Value = 562
$ws.Range("D13").Value = 815
$ws.Range("E13").Value = 682
$ws.Range("G13").Value = 1262
$ws.Range("H13").Value = 1635
